$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "test3" entry (header for column C) is replaced with "日期" (date).
$ws.Range("C1").Value = "日期"

# Move the saved selection/active cell to F9 (was J8).
$ws.Range("F9").Select() | Out-Null
